$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, D value (Price), E value (Volume(1h)) - $null means "leave unchanged"
$updates = @(
    @{Row=2;  D="27.970.60";   E="  +0.00%  "}
    @{Row=3;  D="1.854.63";    E="  -0.79%  "}
    @{Row=4;  D="1.004";       E="  +0.27%  "}
    @{Row=5;  D="311.59";      E="  -0.22%  "}
    @{Row=6;  D=$null;         E="  +0.31%  "}
    @{Row=7;  D="0.5094";      E="  +1.97%  "}
    @{Row=8;  D="0.3806";      E="  -0.64%  "}
    @{Row=9;  D="0.08248";     E="  -7.63%  "}
    @{Row=10; D="1.109";       E="  -1.04%  "}
    @{Row=11; D="41.54";       E="  +0.11%  "}
    @{Row=12; D="6.189";       E="  -2.85%  "}
    @{Row=13; D="20.50";       E="  -1.15%  "}
    @{Row=14; D="1.860.41";    E="  +0.91%  "}
    @{Row=15; D="7.186";       E="  -0.72%  "}
    @{Row=16; D="1.004";       E="  +0.24%  "}
    @{Row=17; D="0.00001096";  E="  -0.34%  "}
    @{Row=18; D="90.52";       E="  -0.69%  "}
    @{Row=19; D=$null;         E="  -1.05%  "}
    @{Row=20; D="17.63";       E="  -2.16%  "}
    @{Row=21; D=$null;         E="  +0.18%  "}
    @{Row=22; D="6.012";       E="  -1.76%  "}
    @{Row=23; D="27.986.18";   E="  -0.08%  "}
    @{Row=24; D="11.01";       E="  -4.30%  "}
    @{Row=25; D="2.242";       E="  -1.84%  "}
    @{Row=26; D="2.539";       E="  +1.15%  "}
    @{Row=27; D="2.071.10";    E="  -0.32%  "}
    @{Row=28; D="157.84";      E="  -0.27%  "}
    @{Row=29; D="20.37";       E="  -1.57%  "}
    @{Row=30; D="124.25";      E="  -1.62%  "}
    @{Row=31; D="0.1055";      E="  -0.51%  "}
    @{Row=32; D="1.035";       E="  -1.98%  "}
    @{Row=33; D="5.611";       E="  +0.35%  "}
    @{Row=34; D="3.598";       E=$null}
    @{Row=35; D="9.430";       E="  +0.35%  "}
    @{Row=36; D="0.06508";     E="  -0.80%  "}
    @{Row=37; D="0.02405";     E="  -0.14%  "}
    @{Row=38; D="0.2162";      E="  -1.32%  "}
    @{Row=39; D="1.201";       E="  +0.01%  "}
    @{Row=40; D="0.6463";      E="  +1.34%  "}
    @{Row=41; D="1.227";       E="  -4.24%  "}
    @{Row=42; D="4.865";       E="  -0.95%  "}
    @{Row=43; D="11.12";       E="  -3.75%  "}
    @{Row=44; D="0.6071";      E="  +1.05%  "}
    @{Row=45; D="13.11";       E="  -1.14%  "}
    @{Row=46; D="1.278";       E="  -0.19%  "}
    @{Row=47; D="3.659";       E="  -0.39%  "}
    @{Row=48; D="1.994";       E="  -0.10%  "}
    @{Row=49; D="1.205";       E="  -1.49%  "}
    @{Row=50; D="119.85";      E="  -0.77%  "}
    @{Row=51; D="78.63";       E="  -0.19%  "}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force text storage so numeric-looking strings (e.g. "311.59")
        # are not coerced into floating point numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
